# "cleaned defensive actions data"
#
# The original sheet had a two-row merged header (row 1 held merged category
# labels like "Tackles"/"Challenges"/"Blocks" over several columns, row 2
# held the real per-column labels) followed directly by the data rows.
#
# The cleaned version turns row 1 into a single, fully-populated header row
# (one distinct label per column, re-using "Blocks"/"Cha"(llenges) plus two
# brand-new labels "Player ID" and "90s"), and keeps the old two-row header
# + a now-empty row 3 around underneath, but hidden, as leftover artifacts
# of the export. The summary "14 Players" row at the bottom is hidden too.
# A couple of data cells that were simply blank (Tkl% for players with no
# attempted tackles) get an explicit 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row-1 category headers spanned merged ranges; unmerge them so
# every column in row 1 can carry its own label.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Populate the new row-1 header labels across A1:W1.
$headerValues = @(
    "Player ID", "Player", "#", "Nation", "Pos", "Age", "90s",
    "Tkl", "TklW", "Def 3rd", "Mid 3rd", "Att 3rd",
    "Cha", "Att", "Tkl%", "Lost",
    "Blocks", "Sh", "Pass", "Int", "Tkl+Int", "Clr", "Err"
)
$headerRow = New-Object 'object[,]' 1,23
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $headerRow[0, $i] = $headerValues[$i]
}
$ws.Range("A1:W1").Value = $headerRow

# Fill in the previously-blank Tkl% cells with an explicit 0.
$ws.Range("O4").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("O17").Value = 0

# The old header row 2, the now-blank row 3, and the trailing summary row 18
# stay in the sheet but are hidden.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(18).Hidden = $true

# Match the saved selection from the source file.
$selectResult = $ws.Range("O19").Select()
